$d = $word.ActiveDocument

# --- Letter addressee block -------------------------------------------------
$d.Paragraphs(6).Range.Text = "Hon. Remar C. Delatina"
$d.Paragraphs(7).Range.Text = "Mayor"
$d.Paragraphs(8).Range.Text = "Pampanga"
$d.Paragraphs(11).Range.Text = "Dear Remar C. Delatina,"

# --- Body paragraph: FY 2020 -> FY 2019 (in-place substring replace) -------
$d.Paragraphs(13).Range.Find.Execute("FY 2020", $false, $false, $false, $false, $false, $true, 1, $false, "FY 2019", 2)

# --- Table: "Hello 1" -> "hello1" --------------------------------------------
$d.Paragraphs(19).Range.Text = "hello1"

# --- Table: drop rows "2"/"Hello 2" and "3"/"Hello 3" -----------------------
$t = $d.Tables(1)
$t.Rows(4).Delete()
$t.Rows(3).Delete()

# --- Signature block ---------------------------------------------------------
$count = $d.Paragraphs.Count
$d.Paragraphs($count - 1).Range.Text = "                                                                         Juan Dela Tina "
$d.Paragraphs($count).Range.Text = "                                                        Regional Director"
